$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (46075 -> 46076, i.e. 2026-02-22 -> 2026-02-23) for every data row (2-11).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 11 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
